# Apply updated crypto price/volume data per diff (rows 2-51, row 29 unchanged)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.528.69"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "3.097.90"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "384.61"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.96"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.93"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "3.586.67"
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.65"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.84"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "3.091.96"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.16"
$ws.Range("E17").Value = "  +8.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.994"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "51.510.28"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  +8.76%  "
$ws.Range("D21").Value = "0.0₃0965"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.37"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.94"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.88"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.13"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.11"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.02"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.32"
$ws.Range("E28").Value = "  -1.68%  "
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.32"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.37"
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0470"
$ws.Range("E34").Value = "  +3.42%  "
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.29"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +2.53%  "
$ws.Range("E39").Value = "  +5.94%  "
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "128.77"
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.56"
$ws.Range("E43").Value = "  -4.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("E44").Value = "  -2.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.49"
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("E47").Value = "  +3.50%  "
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").Value = "2.054.68"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("E51").Value = "  +13.19%  "
